$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.261.21"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.785.85"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.73"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3788"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.73"
$ws.Range("E8").Value = "  -3.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3425"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07490"
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.88"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.471"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "1.789.24"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.088"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06648"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.83"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.636"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.35"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").Value = "27.257.28"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.36"
$ws.Range("E24").Value = "  -6.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.413"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.507"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.543"
$ws.Range("E27").Value = "  -6.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.33"
$ws.Range("E28").Value = "  -4.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.94"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").Value = "1.990.18"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.14"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.000"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.090"
$ws.Range("E33").Value = "  -5.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08699"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.25"
$ws.Range("E35").Value = "  -4.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.667"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6952"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.448"
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2206"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.824"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06330"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02336"
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.44"
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.836"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.151"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.42"
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07131"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.12"
$ws.Range("E51").Value = "  -2.28%  "
